$d = $word.ActiveDocument

$pairs = @(
    @("31×55=1705", "72×33=2376"),
    @("48×61=2928", "50×60=3000"),
    @("84×65=5460", "17×62=1054"),
    @("27×49=1323", "49×70=3430"),
    @("30×32=960", "47×48=2256"),
    @("75×19=1425", "85×42=3570"),
    @("41×70=2870", "12×74=888"),
    @("71×73=5183", "17×61=1037"),
    @("71×55=3905", "18×91=1638"),
    @("27×99=2673", "31×78=2418"),
    @("16×82=1312", "41×68=2788"),
    @("74×45=3330", "38×84=3192"),
    @("11×13=143", "75×83=6225"),
    @("54×60=3240", "34×48=1632"),
    @("50×35=1750", "73×35=2555"),
    @("77×31=2387", "97×53=5141"),
    @("54×31=1674", "71×44=3124"),
    @("56×43=2408", "16×74=1184"),
    @("41×86=3526", "56×65=3640"),
    @("60×40=2400", "11×78=858"),
    @("68×69=4692", "80×79=6320"),
    @("50×23=1150", "89×63=5607"),
    @("14×94=1316", "41×19=779"),
    @("59×35=2065", "69×73=5037"),
    @("67×80=5360", "59×61=3599")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
